$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 'india'
$ws.Range("C62").Value = 'isl'
$ws.Range("D62").Value = '2023-2024'
$ws.Range("E62").Value = 45283.64583333334
$ws.Range("F62").Value = 'Mohun Bagan'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 'Goa'
$ws.Range("I62").Value = 4
$ws.Range("J62").Value = 1.89
$ws.Range("K62").Value = '22/12/2023 03:12'
$ws.Range("L62").Value = 2.55
$ws.Range("M62").Value = '23/12/2023 15:09'
$ws.Range("N62").Value = 3.36
$ws.Range("O62").Value = '22/12/2023 03:12'
$ws.Range("P62").Value = 3.36
$ws.Range("Q62").Value = '23/12/2023 14:59'
$ws.Range("R62").Value = 3.99
$ws.Range("S62").Value = '22/12/2023 03:12'
$ws.Range("T62").Value = 2.83
$ws.Range("U62").Value = '23/12/2023 14:58'
$ws.Range("V62").Value = 'https://www.betexplorer.com/football/india/isl/mohun-bagan-fc-goa/6NvOs69n/'
$ws.Range("A61").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("E61").Copy()
$ws.Range("E62").PasteSpecial(-4122)

# Row 63
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 'india'
$ws.Range("C63").Value = 'isl'
$ws.Range("D63").Value = '2023-2024'
$ws.Range("E63").Value = 45284.54166666666
$ws.Range("F63").Value = 'Bengaluru FC'
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 'North East Utd'
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = 1.67
$ws.Range("K63").Value = '17/12/2023 13:12'
$ws.Range("L63").Value = 1.93
$ws.Range("M63").Value = '24/12/2023 12:26'
$ws.Range("N63").Value = 3.87
$ws.Range("O63").Value = '17/12/2023 13:12'
$ws.Range("P63").Value = 3.63
$ws.Range("Q63").Value = '24/12/2023 12:26'
$ws.Range("R63").Value = 4.5
$ws.Range("S63").Value = '17/12/2023 13:12'
$ws.Range("T63").Value = 3.98
$ws.Range("U63").Value = '24/12/2023 12:26'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/india/isl/bengaluru-fc-north-east-united/pEwStQOh/'
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("E62").Copy()
$ws.Range("E63").PasteSpecial(-4122)

# Row 64
$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 'india'
$ws.Range("C64").Value = 'isl'
$ws.Range("D64").Value = '2023-2024'
$ws.Range("E64").Value = 45284.64583333334
$ws.Range("F64").Value = 'Kerala Blasters'
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 'Mumbai City'
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3.15
$ws.Range("K64").Value = '20/12/2023 15:43'
$ws.Range("L64").Value = 2.95
$ws.Range("M64").Value = '24/12/2023 15:24'
$ws.Range("N64").Value = 3.49
$ws.Range("O64").Value = '20/12/2023 15:43'
$ws.Range("P64").Value = 3.27
$ws.Range("Q64").Value = '24/12/2023 15:23'
$ws.Range("R64").Value = 2.14
$ws.Range("S64").Value = '20/12/2023 15:43'
$ws.Range("T64").Value = 2.51
$ws.Range("U64").Value = '24/12/2023 15:24'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/india/isl/kerala-blasters-mumbai-city/hvwWupva/'
$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("E63").Copy()
$ws.Range("E64").PasteSpecial(-4122)

# Row 65
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 'india'
$ws.Range("C65").Value = 'isl'
$ws.Range("D65").Value = '2023-2024'
$ws.Range("E65").Value = 45286.64583333334
$ws.Range("F65").Value = 'Punjab'
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 'Odisha FC'
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = 2.78
$ws.Range("K65").Value = '23/12/2023 18:42'
$ws.Range("L65").Value = 3.2
$ws.Range("M65").Value = '26/12/2023 15:24'
$ws.Range("N65").Value = 3.28
$ws.Range("O65").Value = '23/12/2023 18:42'
$ws.Range("P65").Value = 3.62
$ws.Range("Q65").Value = '26/12/2023 15:24'
$ws.Range("R65").Value = 2.45
$ws.Range("S65").Value = '23/12/2023 18:42'
$ws.Range("T65").Value = 2.19
$ws.Range("U65").Value = '26/12/2023 15:24'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/india/isl/minerva-punjab-odisha-fc/WUyzu4g5/'
$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("E64").Copy()
$ws.Range("E65").PasteSpecial(-4122)

# Row 66
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 'india'
$ws.Range("C66").Value = 'isl'
$ws.Range("D66").Value = '2023-2024'
$ws.Range("E66").Value = 45287.64583333334
$ws.Range("F66").Value = 'Mohun Bagan'
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 'Kerala Blasters'
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 1.86
$ws.Range("K66").Value = '26/12/2023 17:12'
$ws.Range("L66").Value = 2.04
$ws.Range("M66").Value = '27/12/2023 14:50'
$ws.Range("N66").Value = 3.58
$ws.Range("O66").Value = '26/12/2023 17:12'
$ws.Range("P66").Value = 3.42
$ws.Range("Q66").Value = '27/12/2023 14:50'
$ws.Range("R66").Value = 3.82
$ws.Range("S66").Value = '26/12/2023 17:12'
$ws.Range("T66").Value = 3.8
$ws.Range("U66").Value = '27/12/2023 14:50'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/india/isl/mohun-bagan-kerala-blasters/EgRvvO8B/'
$ws.Range("A65").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("E65").Copy()
$ws.Range("E66").PasteSpecial(-4122)

# Row 67
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 'india'
$ws.Range("C67").Value = 'isl'
$ws.Range("D67").Value = '2023-2024'
$ws.Range("E67").Value = 45288.64583333334
$ws.Range("F67").Value = 'Mumbai City'
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 'Chennaiyin'
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1.51
$ws.Range("K67").Value = '24/12/2023 15:42'
$ws.Range("L67").Value = 1.61
$ws.Range("M67").Value = '28/12/2023 14:49'
$ws.Range("N67").Value = 4.31
$ws.Range("O67").Value = '24/12/2023 15:42'
$ws.Range("P67").Value = 4.12
$ws.Range("Q67").Value = '28/12/2023 14:49'
$ws.Range("R67").Value = 5.23
$ws.Range("S67").Value = '24/12/2023 15:42'
$ws.Range("T67").Value = 5.4
$ws.Range("U67").Value = '28/12/2023 14:56'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/india/isl/mumbai-city-chennaiyin-fc/SdVrwrOH/'
$ws.Range("A66").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("E66").Copy()
$ws.Range("E67").PasteSpecial(-4122)

# Row 68
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 'india'
$ws.Range("C68").Value = 'isl'
$ws.Range("D68").Value = '2023-2024'
$ws.Range("E68").Value = 45289.54166666666
$ws.Range("F68").Value = 'Odisha FC'
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 'Jamshedpur'
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 1.88
$ws.Range("K68").Value = '26/12/2023 15:42'
$ws.Range("L68").Value = 1.93
$ws.Range("M68").Value = '29/12/2023 13:00'
$ws.Range("N68").Value = 3.59
$ws.Range("O68").Value = '26/12/2023 15:42'
$ws.Range("P68").Value = 3.8
$ws.Range("Q68").Value = '29/12/2023 13:00'
$ws.Range("R68").Value = 3.79
$ws.Range("S68").Value = '26/12/2023 15:42'
$ws.Range("T68").Value = 3.78
$ws.Range("U68").Value = '29/12/2023 13:00'
$ws.Range("V68").Value = 'https://www.betexplorer.com/football/india/isl/odisha-fc-jamshedpur/MJUnx2wO/'
$ws.Range("A67").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("E67").Copy()
$ws.Range("E68").PasteSpecial(-4122)

# Row 69
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 'india'
$ws.Range("C69").Value = 'isl'
$ws.Range("D69").Value = '2023-2024'
$ws.Range("E69").Value = 45289.64583333334
$ws.Range("F69").Value = 'North East Utd'
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 'Goa'
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 3.9
$ws.Range("K69").Value = '24/12/2023 13:11'
$ws.Range("L69").Value = 4.74
$ws.Range("M69").Value = '29/12/2023 15:28'
$ws.Range("N69").Value = 3.87
$ws.Range("O69").Value = '24/12/2023 13:11'
$ws.Range("P69").Value = 3.97
$ws.Range("Q69").Value = '29/12/2023 15:29'
$ws.Range("R69").Value = 1.78
$ws.Range("S69").Value = '24/12/2023 13:11'
$ws.Range("T69").Value = 1.71
$ws.Range("U69").Value = '29/12/2023 15:29'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/india/isl/north-east-united-fc-goa/0ATjyMhU/'
$ws.Range("A68").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("E68").Copy()
$ws.Range("E69").PasteSpecial(-4122)

$excel.CutCopyMode = $false